# feat: add new fields to inquiring report (#6)
#
# Insert three new header columns (Sede, Postal Address, Sales TE EMail)
# into the "Data" sheet's header row. "Sede Reason" (previously column O)
# shifts right into P to make room, giving the final header layout:
#   ... N: Domain Reason | O: Sede | P: Sede Reason | Q: Postal Address | R: Sales TE EMail

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Shift "Sede Reason" (old column O) one column to the right, into P,
# freeing up column O for the new "Sede" header.
$ws.Range("O1").Cut($ws.Range("P1"))

# Populate the new header cells.
$ws.Range("O1").Value = "Sede"
$ws.Range("Q1").Value = "Postal Address"
$ws.Range("R1").Value = "Sales TE EMail"

# Match the formatting (gray header fill, etc.) used by the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("O1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column width adjustments.
$ws.Columns.Item(2).ColumnWidth = 21.083333333333332
$ws.Columns.Item(3).ColumnWidth = 18.75
$ws.Columns.Item(15).ColumnWidth = 14.416666666666666
$ws.Columns.Item(16).ColumnWidth = 14.416666666666666
$ws.Columns.Item(17).ColumnWidth = 14.416666666666666
$ws.Columns.Item(18).ColumnWidth = 16.083333333333332

# Extend the AutoFilter to cover the new columns (toggle off first since
# Excel treats re-applying AutoFilter on an overlapping range as a no-op).
$ws.AutoFilterMode = $false
$ws.Range("A1:R1").AutoFilter()

# The hidden _FilterDatabase defined name tracks the AutoFilter range but
# isn't refreshed automatically here, so update it explicitly.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$R`$1"
    }
}

# Update the view: scroll right so column H is the leftmost visible column,
# and move the selection to the cell that was being edited (N14).
$ws.Range("N14").Select()
$excel.ActiveWindow.ScrollColumn = 8
